$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the "+SwitchState(GameState)" paragraph and the
#    "PlayerInput: ..." paragraph (their original positions).
# ------------------------------------------------------------------
$switchPara = $null
$playerInputPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "+SwitchState(GameState)") {
        $switchPara = $i
    }
    if ($t -eq "PlayerInput: post events by delegate, ShipMove, PlayerLoot, ShipShot will receive it") {
        $playerInputPara = $i
    }
}

# ------------------------------------------------------------------
# 2) Remove the old "PlayerInput..." paragraph (it will be re-created
#    right after "+SwitchState(GameState)").
# ------------------------------------------------------------------
$oldP = $d.Paragraphs.Item($playerInputPara)
$oldP.Range.Delete()

# ------------------------------------------------------------------
# 3) Turn "+SwitchState(GameState)" into a paragraph whose mark is
#    bold, and append a trailing bold space run to it.
# ------------------------------------------------------------------
$switchP = $d.Paragraphs.Item($switchPara)
$switchR = $switchP.Range
$switchXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>+SwitchState(GameState)</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$switchR.InsertXML($switchXml)

# ------------------------------------------------------------------
# 4) Insert a brand-new paragraph right after it holding the
#    "PlayerInput: post events by delegate, ..." text (4 runs, the
#    ":" and following space merged into a single run this time).
# ------------------------------------------------------------------
$switchP = $d.Paragraphs.Item($switchPara)
$switchP.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($switchPara + 1)
$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>PlayerInput</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>post events by delegate</w:t></w:r><w:r><w:t>, ShipMove, PlayerLoot, ShipShot will receive it</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newXml)

# ------------------------------------------------------------------
# 5) Insert one more, genuinely empty paragraph right after the
#    PlayerInput paragraph (matching the lone "<w:p/>" added by the
#    diff). Build it via InsertXML rather than InsertParagraphAfter
#    so it does not inherit the preceding bold run formatting.
# ------------------------------------------------------------------
$newPara = $d.Paragraphs.Item($switchPara + 1)
$newPara.Range.InsertParagraphAfter()
$emptyPara = $d.Paragraphs.Item($switchPara + 2)
$emptyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$emptyPara.Range.InsertXML($emptyXml)
